# Update "想去人数" (F column) counts that changed between site crawls.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows, so the
# update has to be applied in both places (F4 diverges slightly between the
# two sheets: 2942 vs 2943).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 189
    4  = 2942
    5  = 210
    16 = 223
    18 = 225
    23 = 346
    24 = 129
    27 = 1951
    28 = 49
    31 = 155
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates2 = @{
    3  = 189
    4  = 2943
    5  = 210
    16 = 223
    18 = 225
    23 = 346
    24 = 129
    27 = 1951
    28 = 49
    31 = 155
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates2.Keys) {
    $ws4.Range("F$row").Value = $updates2[$row]
}
